$d = $word.ActiveDocument

# --- Edit 1: first paragraph ("This is a Microsoft word document.") ---
# The original single run keeps its text but gains two trailing spaces;
# three new runs colored red (FF0000) are appended after it, before the
# paragraph mark, spelling out "(This is a change - Version for main
# branch)" (with a real en dash) split across three runs.
$p1 = $d.Paragraphs(1)
$newFirstText = "This is a Microsoft word document.  "
$origLen = $p1.Range.Text.Length - 1   # exclude the paragraph mark
$rng = $d.Range($p1.Range.Start, $p1.Range.Start + $origLen)
$rng.Text = $newFirstText

$pos = $p1.Range.Start + $newFirstText.Length

$chunk1 = "(This is a change " + [char]0x2013 + " Ve"
$chunk2 = "rsion for main branch"
$chunk3 = ")"

$ip1 = $d.Range($pos, $pos)
$ip1.InsertAfter($chunk1)
$d.Range($pos, $pos + $chunk1.Length).Font.Color = 255
$pos = $pos + $chunk1.Length

$ip2 = $d.Range($pos, $pos)
$ip2.InsertAfter($chunk2)
$d.Range($pos, $pos + $chunk2.Length).Font.Color = 255
$pos = $pos + $chunk2.Length

$ip3 = $d.Range($pos, $pos)
$ip3.InsertAfter($chunk3)
$d.Range($pos, $pos + $chunk3.Length).Font.Color = 255
$pos = $pos + $chunk3.Length

# --- Edit 2: the empty paragraph right after "It will be treated as a
# binary file by Git." (paraId 7476926D) currently carries Menlo/black/
# 9pt mark formatting plus w:textAlignment="baseline" in its pPr. It
# becomes a completely bare paragraph, i.e. <w:p/>. ---
$p3 = $d.Paragraphs(3)
$null = $p3.Range.InsertXML("<w:p/>")

Write-Output "done"
